# Merge the three separately-formatted runs that make up each
#   <id>p124v_N</id>
# paragraph ("<id>", "p124v_N", "</id>") into a single run, using the
# formatting of the surrounding "<id>"/"</id>" runs (Courier New, 18pt,
# color 7f6000). Paragraphs such as "<id>fig_p124v_N</id>" are left
# untouched.

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range.Text

    # Strip the trailing paragraph mark before comparing.
    $idText = $full.TrimEnd([char]13)

    if ($idText -match "^<id>p124v_\d+</id>$") {
        $pStart = $p.Range.Start
        $len = $idText.Length
        $pEnd = $pStart + $len

        # Force Word to collapse the run(s) covering the text into a single
        # run by assigning FormattedText with different placeholder content
        # (same length); the resulting single run inherits the formatting
        # of the range's first original run ("<id>"/"</id>" - Courier New,
        # 18pt, color 7f6000).
        $target = $d.Range($pStart, $pEnd)
        $placeholder = $target.FormattedText
        $placeholder.Text = "".PadRight($len, 'P')
        $target.FormattedText = $placeholder

        # Now restore the real text within that now-uniform run.
        $fixed = $d.Range($pStart, $pStart + $len)
        $fixed.Text = $idText
    }
}
